$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C cell values to append power_set info (experiment 1 code / colorado gas model)
$ws.Range("C2").Value = "(5000, 0), p:(5100, 2), t:(5200, 4), power_set: (6000)"
$ws.Range("C3").Value = "(5001, 0), p:(5101, 2), t:(5201, 4), 6001"
$ws.Range("C4").Value = "(5003, 0), p:(5103, 2), t:(5203, 4), 6003"
$ws.Range("C5").Value = "(5004, 0), p:(5104, 2), t:(5204, 4), 6004"
$ws.Range("C6").Value = "(5005, 0), p:(5105, 2), t:(5205, 4), 6005"
$ws.Range("C7").Value = "(5006, 0), p:(5106, 2), t:(5206, 4), 6006"
$ws.Range("C8").Value = "(5007, 0), p:(5107 , 2), t:(5207, 4), 6007"

# Update column C width (COM ColumnWidth differs from stored OOXML width by ~0.8333;
# 58.67 maps to a stored width of 59.5)
$ws.Columns.Item(3).ColumnWidth = 58.67

# Update the active selection to C9
$ws.Range("C9").Select()
